# Update "paises.xlsx" (countries & provincias Spain) data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 16:43"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 3549519
$ws.Range("C4").Value = 4442
$ws.Range("D4").Value = 1601030
$ws.Range("E4").Value = 1809287
$ws.Range("G4").Value = 59
$ws.Range("H4").Value = 139202

# --- India (row 6) ---
$ws.Range("B6").Value = 956992
$ws.Range("C6").Value = 19505
$ws.Range("D6").Value = 604624
$ws.Range("E6").Value = 327665
$ws.Range("G6").Value = 388
$ws.Range("H6").Value = 24703

# --- Serbia (row 62) ---
$ws.Range("D62").Value = 14047
$ws.Range("E62").Value = 4858

# --- Kenia (row 73) ---
$ws.Range("D73").Value = 3068
$ws.Range("E73").Value = 7975
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 209

# --- Noruega (row 78) ---
$ws.Range("B78").Value = 9006
$ws.Range("C78").Value = 5
$ws.Range("E78").Value = 615

# --- Tayikistan (row 90) ---
$ws.Range("B90").Value = 6695
$ws.Range("C90").Value = 52
$ws.Range("D90").Value = 5383
$ws.Range("E90").Value = 1256

# --- Sri Lanka (row 110) ---
$ws.Range("B110").Value = 2671
$ws.Range("C110").Value = 6
$ws.Range("D110").Value = 2001
$ws.Range("E110").Value = 659

# --- Georgia (row 144) ---
$ws.Range("B144").Value = 1004
$ws.Range("C144").Value = 5
$ws.Range("E144").Value = 116

# --- Swap Groenlandia / Islas Malvinas labels (rows 209-210) ---
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
